# CTMS: Added test step - assigned questionnaire
# Adds a new "AssignedQuestionnaire" test step (columns G/H) to the
# existing AddVisitScheduleForSite step block on row 4 (header) / row 5 (data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (bold, matching the existing F4 "Category" header style)
$ws.Range("G4").Value = "QuestionnaireAssignedTo"
$ws.Range("G5").Value = "Study"
$ws.Range("H4").Value = "Template Name"
$ws.Range("H5").Value = "SQV Report"

$ws.Range("G4").Font.Bold = $true
$ws.Range("H4").Font.Bold = $true

# Move the active selection to follow the newly populated range, like Excel
# would after tabbing through data entry.
$ws.Range("H8").Select() | Out-Null
